# Update scripts with new TPM-derived values.
# The sheet previously held 3 rows (MuSCs -> ECs / FAPs / MuSCs as
# sending/target clusters for the Fgf22/Fgfr1 pair). The new data adds a
# second sending cluster (ECs) with its own recomputed edge weights, while
# also recomputing the existing MuSCs-sourced rows, and inserts 3 more rows
# (5:7) holding the MuSCs-sourced results below the ECs-sourced ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> Fgf22 -> Fgfr1 -> ECs -----------------------------------
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf22"
$ws.Range("C2").Value = "Fgfr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09830299999999999
$ws.Range("H2").Value = 0.294909
$ws.Range("I2").Value = 0.3722197960868258
$ws.Range("J2").Value = 0.3722197960868259
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.294987
$ws.Range("N2").Value = 6.884961000000001
$ws.Range("O2").Value = 0.0158275801650097
$ws.Range("P2").Value = 0.0158275801650097
$ws.Range("Q2").Value = 0.225604107061
$ws.Range("R2").Value = 2.030436963549
$ws.Range("S2").Value = 0.0058913386615678
$ws.Range("T2").Value = 0.005891338661567802

# --- Row 3: ECs -> Fgf22 -> Fgfr1 -> FAPs -----------------------------------
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf22"
$ws.Range("C3").Value = "Fgfr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09830299999999999
$ws.Range("H3").Value = 0.294909
$ws.Range("I3").Value = 0.3722197960868258
$ws.Range("J3").Value = 0.3722197960868259
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 111.5917106666667
$ws.Range("N3").Value = 334.775132
$ws.Range("O3").Value = 0.769602070219672
$ws.Range("P3").Value = 0.7696020702196722
$ws.Range("Q3").Value = 10.96979993366533
$ws.Range("R3").Value = 98.72819940298798
$ws.Range("S3").Value = 0.2864611256451653
$ws.Range("T3").Value = 0.2864611256451655

# --- Row 4: ECs -> Fgf22 -> Fgfr1 -> MuSCs ----------------------------------
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf22"
$ws.Range("C4").Value = "Fgfr1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09830299999999999
$ws.Range("H4").Value = 0.294909
$ws.Range("I4").Value = 0.3722197960868258
$ws.Range("J4").Value = 0.3722197960868259
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2145703496153182
$ws.Range("P4").Value = 0.2145703496153182
$ws.Range("Q4").Value = 3.058455659175666
$ws.Range("R4").Value = 27.526100932581
$ws.Range("S4").Value = 0.07986733178009267
$ws.Range("T4").Value = 0.0798673317800927

# --- Insert 3 fresh rows at 5:7 for the MuSCs-sourced results --------------
$ws.Range("A5:A7").EntireRow.Insert()

# --- Row 5: MuSCs -> Fgf22 -> Fgfr1 -> ECs ----------------------------------
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf22"
$ws.Range("C5").Value = "Fgfr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1657963333333334
$ws.Range("H5").Value = 0.497389
$ws.Range("I5").Value = 0.6277802039131741
$ws.Range("J5").Value = 0.6277802039131741
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.294987
$ws.Range("N5").Value = 6.884961000000001
$ws.Range("O5").Value = 0.0158275801650097
$ws.Range("P5").Value = 0.0158275801650097
$ws.Range("Q5").Value = 0.3805004296476668
$ws.Range("R5").Value = 3.424503866829001
$ws.Range("S5").Value = 0.009936241503441901
$ws.Range("T5").Value = 0.009936241503441901

# --- Row 6: MuSCs -> Fgf22 -> Fgfr1 -> FAPs ---------------------------------
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf22"
$ws.Range("C6").Value = "Fgfr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1657963333333334
$ws.Range("H6").Value = 0.497389
$ws.Range("I6").Value = 0.6277802039131741
$ws.Range("J6").Value = 0.6277802039131741
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 111.5917106666667
$ws.Range("N6").Value = 334.775132
$ws.Range("O6").Value = 0.769602070219672
$ws.Range("P6").Value = 0.7696020702196722
$ws.Range("Q6").Value = 18.50149645892756
$ws.Range("R6").Value = 166.513468130348
$ws.Range("S6").Value = 0.4831409445745066
$ws.Range("T6").Value = 0.4831409445745068

# --- Row 7: MuSCs -> Fgf22 -> Fgfr1 -> MuSCs --------------------------------
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf22"
$ws.Range("C7").Value = "Fgfr1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1657963333333334
$ws.Range("H7").Value = 0.497389
$ws.Range("I7").Value = 0.6277802039131741
$ws.Range("J7").Value = 0.6277802039131741
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 31.11253633333333
$ws.Range("N7").Value = 93.337609
$ws.Range("O7").Value = 0.2145703496153182
$ws.Range("P7").Value = 0.2145703496153182
$ws.Range("Q7").Value = 5.158344444766779
$ws.Range("R7").Value = 46.425100002901
$ws.Range("S7").Value = 0.1347030178352255
$ws.Range("T7").Value = 0.1347030178352255
